# Auto-generated edit script: updates Anima_Profits market-price derived columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled-runner refresh diff.

$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 19
$ws_ALC.Range("H19").Value = 6661.9697
$ws_ALC.Range("I19").Value = 638.6667
$ws_ALC.Range("J19").Value = 8920.708000000001
$ws_ALC.Range("K19").Value = 638.6667
$ws_ALC.Range("L19").Value = 8920.708000000001
$ws_ALC.Range("M19").Value = -463.6667
$ws_ALC.Range("N19").Value = -9270.708000000001

# ALC row 96
$ws_ALC.Range("H96").Value = 1222.5
$ws_ALC.Range("I96").Value = 733.25
$ws_ALC.Range("J96").Value = 1711.75
$ws_ALC.Range("K96").Value = 2199.75
$ws_ALC.Range("L96").Value = 5135.25
$ws_ALC.Range("M96").Value = -826.75
$ws_ALC.Range("N96").Value = -7881.25

# ALC row 100
$ws_ALC.Range("H100").Value = 2980.4
$ws_ALC.Range("I100").Value = 2985
$ws_ALC.Range("K100").Value = 2985
$ws_ALC.Range("M100").Value = -2444

# ALC row 137
$ws_ALC.Range("H137").Value = 3139.6
$ws_ALC.Range("I137").Value = 2188.4
$ws_ALC.Range("J137").Value = 9798
$ws_ALC.Range("K137").Value = 6565.200000000001
$ws_ALC.Range("L137").Value = 29394
$ws_ALC.Range("M137").Value = -4015.200000000001
$ws_ALC.Range("N137").Value = -34494

# ALC row 138
$ws_ALC.Range("H138").Value = 2079.11
$ws_ALC.Range("I138").Value = 690.8182
$ws_ALC.Range("J138").Value = 2250.6965
$ws_ALC.Range("K138").Value = 2072.4546
$ws_ALC.Range("L138").Value = 6752.0895
$ws_ALC.Range("M138").Value = 3067.5454
$ws_ALC.Range("N138").Value = -17032.0895

# ARM row 32
$ws_ARM.Range("H32").Value = 381180.4
$ws_ARM.Range("I32").Value = 423089.62
$ws_ARM.Range("J32").Value = 22623.777
$ws_ARM.Range("K32").Value = 423089.62
$ws_ARM.Range("L32").Value = 22623.777
$ws_ARM.Range("M32").Value = -422802.62
$ws_ARM.Range("N32").Value = -23197.777

# ARM row 61
$ws_ARM.Range("H61").Value = 2751.4614
$ws_ARM.Range("I61").Value = 2274.24
$ws_ARM.Range("J61").Value = 3603.6428
$ws_ARM.Range("K61").Value = 2274.24
$ws_ARM.Range("L61").Value = 3603.6428
$ws_ARM.Range("M61").Value = -2062.24
$ws_ARM.Range("N61").Value = -4027.6428

# ARM row 102
$ws_ARM.Range("H102").Value = 1951.3684
$ws_ARM.Range("I102").Value = 1922.1177
$ws_ARM.Range("K102").Value = 1922.1177
$ws_ARM.Range("M102").Value = -300.1177

# ARM row 132
$ws_ARM.Range("H132").Value = 4454.7295
$ws_ARM.Range("I132").Value = 3340.7896
$ws_ARM.Range("J132").Value = 5630.5557
$ws_ARM.Range("K132").Value = 10022.3688
$ws_ARM.Range("L132").Value = 16891.6671
$ws_ARM.Range("M132").Value = -7492.3688
$ws_ARM.Range("N132").Value = -21951.6671

# ARM row 136
$ws_ARM.Range("H136").Value = 2751.4614
$ws_ARM.Range("I136").Value = 2274.24
$ws_ARM.Range("J136").Value = 3603.6428
$ws_ARM.Range("K136").Value = 6822.719999999999
$ws_ARM.Range("L136").Value = 10810.9284
$ws_ARM.Range("M136").Value = -4272.719999999999
$ws_ARM.Range("N136").Value = -15910.9284

# BSM row 94
$ws_BSM.Range("H94").Value = 1825
$ws_BSM.Range("I94").Value = 1450
$ws_BSM.Range("J94").Value = 1900
$ws_BSM.Range("K94").Value = 1450
$ws_BSM.Range("L94").Value = 1900
$ws_BSM.Range("M94").Value = -999
$ws_BSM.Range("N94").Value = -2802

# BSM row 105
$ws_BSM.Range("H105").Value = 13892064
$ws_BSM.Range("I105").Value = 17860212
$ws_BSM.Range("J105").Value = 3550
$ws_BSM.Range("K105").Value = 17860212
$ws_BSM.Range("L105").Value = 3550
$ws_BSM.Range("M105").Value = -17858465
$ws_BSM.Range("N105").Value = -7044

# BSM row 134
$ws_BSM.Range("H134").Value = 2354.2432
$ws_BSM.Range("I134").Value = 2140.682
$ws_BSM.Range("J134").Value = 2667.4666
$ws_BSM.Range("K134").Value = 6422.045999999999
$ws_BSM.Range("L134").Value = 8002.399800000001
$ws_BSM.Range("M134").Value = -3887.045999999999
$ws_BSM.Range("N134").Value = -13072.3998

# CRP row 31
$ws_CRP.Range("H31").Value = 8258.206
$ws_CRP.Range("I31").Value = 2308.8462
$ws_CRP.Range("J31").Value = 11941.143
$ws_CRP.Range("K31").Value = 2308.8462
$ws_CRP.Range("L31").Value = 11941.143
$ws_CRP.Range("M31").Value = -2013.8462
$ws_CRP.Range("N31").Value = -12531.143

# CRP row 34
$ws_CRP.Range("H34").Value = 8258.206
$ws_CRP.Range("I34").Value = 2308.8462
$ws_CRP.Range("J34").Value = 11941.143
$ws_CRP.Range("K34").Value = 2308.8462
$ws_CRP.Range("L34").Value = 11941.143
$ws_CRP.Range("M34").Value = -2106.8462
$ws_CRP.Range("N34").Value = -12345.143

# CRP row 58
$ws_CRP.Range("H58").Value = 1327.3793
$ws_CRP.Range("I58").Value = 1019.125
$ws_CRP.Range("K58").Value = 1019.125
$ws_CRP.Range("M58").Value = -816.125

# CRP row 105
$ws_CRP.Range("H105").Value = 1933.5
$ws_CRP.Range("J105").Value = 1833.3334
$ws_CRP.Range("L105").Value = 1833.3334
$ws_CRP.Range("N105").Value = -5327.3334

# CRP row 132
$ws_CRP.Range("H132").Value = 7248620
$ws_CRP.Range("I132").Value = 1858.4375
$ws_CRP.Range("J132").Value = 23812646
$ws_CRP.Range("K132").Value = 5575.3125
$ws_CRP.Range("L132").Value = 71437938
$ws_CRP.Range("M132").Value = -3045.3125
$ws_CRP.Range("N132").Value = -71442998

# CRP row 134
$ws_CRP.Range("H134").Value = 1617.9259
$ws_CRP.Range("I134").Value = 1459.5
$ws_CRP.Range("K134").Value = 4378.5
$ws_CRP.Range("M134").Value = -1843.5

# CRP row 136
$ws_CRP.Range("H136").Value = 1327.3793
$ws_CRP.Range("I136").Value = 1019.125
$ws_CRP.Range("K136").Value = 3057.375
$ws_CRP.Range("M136").Value = -507.375

# CUL row 33
$ws_CUL.Range("H33").Value = 20140.2
$ws_CUL.Range("I33").Value = 33433
$ws_CUL.Range("J33").Value = 201
$ws_CUL.Range("K33").Value = 200598
$ws_CUL.Range("L33").Value = 1206
$ws_CUL.Range("M33").Value = -200315
$ws_CUL.Range("N33").Value = -1772

# CUL row 114
$ws_CUL.Range("H114").Value = 1050.05
$ws_CUL.Range("I114").Value = 166.91667
$ws_CUL.Range("J114").Value = 2374.75
$ws_CUL.Range("K114").Value = 500.75001
$ws_CUL.Range("L114").Value = 7124.25
$ws_CUL.Range("M114").Value = 2753.24999
$ws_CUL.Range("N114").Value = -13632.25

# GSM row 102
$ws_GSM.Range("H102").Value = 1502.4
$ws_GSM.Range("I102").Value = 1378
$ws_GSM.Range("J102").Value = 2000
$ws_GSM.Range("K102").Value = 1378
$ws_GSM.Range("L102").Value = 2000
$ws_GSM.Range("M102").Value = 244
$ws_GSM.Range("N102").Value = -5244

# GSM row 126
$ws_GSM.Range("H126").Value = 1982.4615
$ws_GSM.Range("I126").Value = 1982.4615
$ws_GSM.Range("J126").Value = 0
$ws_GSM.Range("K126").Value = 5947.3845
$ws_GSM.Range("L126").Value = 0
$ws_GSM.Range("M126").Value = -3477.3845
$ws_GSM.Range("N126").ClearContents()

# GSM row 132
$ws_GSM.Range("H132").Value = 2516.8462
$ws_GSM.Range("I132").Value = 2455
$ws_GSM.Range("J132").Value = 2562.2
$ws_GSM.Range("K132").Value = 7365
$ws_GSM.Range("L132").Value = 7686.599999999999
$ws_GSM.Range("M132").Value = -4835
$ws_GSM.Range("N132").Value = -12746.6

# LTW row 132
$ws_LTW.Range("H132").Value = 2373.8
$ws_LTW.Range("I132").Value = 1408.4546
$ws_LTW.Range("K132").Value = 4225.3638
$ws_LTW.Range("M132").Value = -1695.3638

# WVR row 46
$ws_WVR.Range("H46").Value = 98429
$ws_WVR.Range("J46").Value = 98429
$ws_WVR.Range("L46").Value = 98429
$ws_WVR.Range("N46").Value = -98891

# WVR row 132
$ws_WVR.Range("H132").Value = 3625199.5
$ws_WVR.Range("I132").Value = 1851.2759
$ws_WVR.Range("K132").Value = 5553.8277
$ws_WVR.Range("M132").Value = -3023.8277

# WVR row 134
$ws_WVR.Range("H134").Value = 98429
$ws_WVR.Range("J134").Value = 98429
$ws_WVR.Range("L134").Value = 295287
$ws_WVR.Range("N134").Value = -300357

# WVR row 136
$ws_WVR.Range("H136").Value = 2428.1282
$ws_WVR.Range("I136").Value = 2268.1538
$ws_WVR.Range("K136").Value = 6804.4614
$ws_WVR.Range("M136").Value = -4254.4614
